$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: "Tabla perforata vopsita electrostatic" -> new, more specific note
$ws.Range("D2").Value = "Tabla perforata vopsita electrostatic. De la 1 la 4 "

# New G25:G27 cells (3D asset column) - copy formatting from the sibling
# F-column cell on the same row so the style matches the rest of the row,
# then set the value.
$ws.Range("F25").Copy($ws.Range("G25"))
$ws.Range("F26").Copy($ws.Range("G26"))
$ws.Range("F27").Copy($ws.Range("G27"))

$ws.Range("G25").Value = "assets/model17.glb"
$ws.Range("G26").Value = "assets/model17.glb"
$ws.Range("G27").Value = "assets/model17.glb"

# Update the on-screen selection to match the saved view state.
$ws.Range("I25").Select()
